$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20 (ALC)
$ws.Range("H20").Value = 5704.2
$ws.Range("I20").Value = 3380.25
$ws.Range("J20").Value = 15000
$ws.Range("K20").Value = 3380.25
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = -3150.25
$ws.Range("N20").Value = -15460

# Row 28 (ALC)
$ws.Range("H28").Value = 1063
$ws.Range("I28").Value = 894.5
$ws.Range("J28").Value = 1400
$ws.Range("K28").Value = 894.5
$ws.Range("L28").Value = 1400
$ws.Range("M28").Value = -409.5
$ws.Range("N28").Value = -2370

# Row 35 (ALC)
$ws.Range("H35").Value = 5704.2
$ws.Range("I35").Value = 3380.25
$ws.Range("J35").Value = 15000
$ws.Range("K35").Value = 3380.25
$ws.Range("L35").Value = 15000
$ws.Range("M35").Value = -3001.25
$ws.Range("N35").Value = -15758

# Row 70 (ALC)
$ws.Range("H70").Value = 1250
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -5040

# Row 73 (ALC)
$ws.Range("H73").Value = 1250
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -6372

# Row 94 (ALC)
$ws.Range("H94").Value = 4565.8335
$ws.Range("I94").Value = 4279
$ws.Range("J94").Value = 6000
$ws.Range("K94").Value = 4279
$ws.Range("L94").Value = 6000
$ws.Range("M94").Value = -3828

# Row 100 (ALC)
$ws.Range("H100").Value = 1932
$ws.Range("I100").Value = 1932
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1932
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1391

# Row 107 (ALC)
$ws.Range("H107").Value = 191
$ws.Range("I107").Value = 191
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 191
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1729

# Row 138 (ALC)
$ws.Range("H138").Value = 3706493.5
$ws.Range("I138").Value = 20001220
$ws.Range("J138").Value = 3146.7273
$ws.Range("K138").Value = 60003660
$ws.Range("L138").Value = 9440.1819
$ws.Range("M138").Value = -59998520
$ws.Range("N138").Value = -19720.1819

$ws = $wb.Worksheets.Item("ARM")
# Row 39 (ARM)
$ws.Range("H39").Value = 3500
$ws.Range("I39").Value = 3500
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 3500
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -2980

# Row 45 (ARM)
$ws.Range("H45").Value = 3584.1333
$ws.Range("I45").Value = 3248.9167
$ws.Range("J45").Value = 4925
$ws.Range("K45").Value = 3248.9167
$ws.Range("L45").Value = 4925
$ws.Range("M45").Value = -2871.9167
$ws.Range("N45").Value = -5679

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (BSM)
$ws.Range("H20").Value = 6000
$ws.Range("I20").Value = 4000
$ws.Range("J20").Value = 8000
$ws.Range("K20").Value = 4000
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = -3753
$ws.Range("N20").Value = -8494

# Row 38 (BSM)
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 10 (CRP)
$ws.Range("H10").Value = 521.8333
$ws.Range("I10").Value = 165.28572
$ws.Range("J10").Value = 1021
$ws.Range("K10").Value = 165.28572
$ws.Range("L10").Value = 1021
$ws.Range("M10").Value = -26.28572
$ws.Range("N10").Value = -1299

$ws = $wb.Worksheets.Item("CUL")
# Row 12 (CUL)
$ws.Range("H12").Value = 39.153847
$ws.Range("I12").Value = 45.75
$ws.Range("J12").Value = 36.22222
$ws.Range("K12").Value = 137.25
$ws.Range("L12").Value = 108.66666
$ws.Range("M12").Value = 35.75
$ws.Range("N12").Value = -454.66666

# Row 22 (CUL)
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()

# Row 27 (CUL)
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()

# Row 46 (CUL)
$ws.Range("H46").Value = 1750
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1750
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5250
$ws.Range("N46").Value = -5432

# Row 81 (CUL)
$ws.Range("H81").Value = 1721
$ws.Range("I81").Value = 1856.5
$ws.Range("J81").Value = 1450
$ws.Range("K81").Value = 5569.5
$ws.Range("L81").Value = 4350
$ws.Range("M81").Value = -4446.5
$ws.Range("N81").Value = -6596

# Row 84 (CUL)
$ws.Range("H84").Value = 1721
$ws.Range("I84").Value = 1856.5
$ws.Range("J84").Value = 1450
$ws.Range("K84").Value = 16708.5
$ws.Range("L84").Value = 13050
$ws.Range("M84").Value = -11092.5
$ws.Range("N84").Value = -24282

$ws = $wb.Worksheets.Item("GSM")
# Row 21 (GSM)
$ws.Range("H21").Value = 25000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 25000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 25000
$ws.Range("N21").Value = -25346

# Row 30 (GSM)
$ws.Range("H30").Value = 25000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 25000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 25000
$ws.Range("N30").Value = -25210

# Row 31 (GSM)
$ws.Range("H31").Value = 1900
$ws.Range("I31").Value = 1900
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1900
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1608

# Row 37 (GSM)
$ws.Range("H37").Value = 1900
$ws.Range("I37").Value = 1900
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1900
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -1623

# Row 80 (GSM)
$ws.Range("H80").Value = 999
$ws.Range("I80").Value = 999
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 999
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1
$ws.Range("N80").ClearContents()

# Row 83 (GSM)
$ws.Range("H83").Value = 999
$ws.Range("I83").Value = 999
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 4995
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -3
$ws.Range("N83").ClearContents()

# Row 102 (GSM)
$ws.Range("H102").Value = 3137.8462
$ws.Range("I102").Value = 1724.5
$ws.Range("J102").Value = 7849
$ws.Range("K102").Value = 1724.5
$ws.Range("L102").Value = 7849
$ws.Range("M102").Value = -102.5
$ws.Range("N102").Value = -11093

# Row 104 (GSM)
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Row 132 (GSM)
$ws.Range("H132").Value = 2198.4285
$ws.Range("I132").Value = 1731.5
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5194.5
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -2664.5
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (LTW)
$ws.Range("H40").Value = 2333.3333
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -1864
$ws.Range("N40").Value = -2772

# Row 55 (LTW)
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

# Row 64 (LTW)
$ws.Range("H64").Value = 9575
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 9575
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 9575
$ws.Range("N64").Value = -10025

# Row 67 (LTW)
$ws.Range("H67").Value = 9575
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 9575
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 9575
$ws.Range("N67").Value = -11135

# Row 93 (LTW)
$ws.Range("H93").Value = 1482.7273
$ws.Range("I93").Value = 1482.7273
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1482.7273
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -234.7273

$ws = $wb.Worksheets.Item("WVR")
# Row 7 (WVR)
$ws.Range("H7").Value = 1900
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 1900
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1900
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -2126

# Row 9 (WVR)
$ws.Range("H9").Value = 6
$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 6
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 134
$ws.Range("N9").ClearContents()

# Row 43 (WVR)
$ws.Range("H43").Value = 17500
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 17500
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 17500
$ws.Range("N43").Value = -17798

# Row 122 (WVR)
$ws.Range("H122").Value = 1786.6364
$ws.Range("I122").Value = 1710.625
$ws.Range("J122").Value = 1989.3334
$ws.Range("K122").Value = 5131.875
$ws.Range("L122").Value = 5968.0002
$ws.Range("M122").Value = -2681.875
$ws.Range("N122").Value = -10868.0002

# Row 140 (WVR)
$ws.Range("H140").Value = 89999.664
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 89999.664
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 89999.664
$ws.Range("N140").Value = -100359.664
